# Auto-generated edit script: applies the cryptos.xlsx price/volume update diff
# (commit: "Updated cryptos list on Sat Jan 27 02:36:44 UTC 2024 with GitHub Actions")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to text format so numeric-looking strings
# (e.g. "41.901.76", "34.10") are preserved verbatim as text, matching
# the original inline-string cell contents instead of being parsed as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value2 = '41.901.76'
$ws.Range("E2").Value2 = '  +4.27%  '
$ws.Range("D3").Value2 = '2.272.87'
$ws.Range("E3").Value2 = '  +1.84%  '
$ws.Range("E4").Value2 = '  -0.05%  '
$ws.Range("D5").Value2 = '304.07'
$ws.Range("E5").Value2 = '  +3.32%  '
$ws.Range("D6").Value2 = '93.44'
$ws.Range("E6").Value2 = '  +6.29%  '
$ws.Range("D7").Value2 = '0.534'
$ws.Range("E7").Value2 = '  +4.05%  '
$ws.Range("E8").Value2 = '  -0.05%  '
$ws.Range("D9").Value2 = '0.488'
$ws.Range("E9").Value2 = '  +3.83%  '
$ws.Range("D10").Value2 = '32.79'
$ws.Range("E10").Value2 = '  +6.93%  '
$ws.Range("D11").Value2 = '53.23'
$ws.Range("E11").Value2 = '  +4.42%  '
$ws.Range("E12").Value2 = '  +2.42%  '
$ws.Range("E13").Value2 = '  +2.67%  '
$ws.Range("E14").Value2 = '  +3.64%  '
$ws.Range("D15").Value2 = '2.624.39'
$ws.Range("E15").Value2 = '  +1.69%  '
$ws.Range("D16").Value2 = '14.29'
$ws.Range("E16").Value2 = '  +3.16%  '
$ws.Range("D17").Value2 = '2.267.77'
$ws.Range("E17").Value2 = '  -0.25%  '
$ws.Range("D18").Value2 = '0.763'
$ws.Range("E18").Value2 = '  +3.51%  '
$ws.Range("D19").Value2 = '41.824.56'
$ws.Range("E19").Value2 = '  +4.25%  '
$ws.Range("D20").Value2 = '12.29'
$ws.Range("E20").Value2 = '  +8.89%  '
$ws.Range("D21").Value2 = '0.0₃0909'
$ws.Range("E21").Value2 = '  +2.31%  '
$ws.Range("D22").Value2 = '5.96'
$ws.Range("E22").Value2 = '  +2.82%  '
$ws.Range("D23").Value2 = '67.37'
$ws.Range("E23").Value2 = '  +2.28%  '
$ws.Range("D24").Value2 = '243.49'
$ws.Range("E24").Value2 = '  +2.68%  '
$ws.Range("D25").Value2 = '2.58'
$ws.Range("E25").Value2 = '  +3.77%  '
$ws.Range("E26").Value2 = '  -0.09%  '
$ws.Range("D27").Value2 = '1.93'
$ws.Range("E27").Value2 = '  +4.78%  '
$ws.Range("D28").Value2 = '24.32'
$ws.Range("D29").Value2 = '9.64'
$ws.Range("E29").Value2 = '  +3.47%  '
$ws.Range("E30").Value2 = '  -3.80%  '
$ws.Range("B31").Value2 = 'InjectiveProtocol'
$ws.Range("C31").Value2 = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D31").Value2 = '34.10'
$ws.Range("E31").Value2 = '  +6.97%  '
$ws.Range("B32").Value2 = 'Monero'
$ws.Range("C32").Value2 = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D32").Value2 = '158.35'
$ws.Range("E32").Value2 = '  +0.38%  '
$ws.Range("E33").Value2 = '  -0.01%  '
$ws.Range("E34").Value2 = '  +4.17%  '
$ws.Range("D35").Value2 = '0.0753'
$ws.Range("E35").Value2 = '  +4.97%  '
$ws.Range("E36").Value2 = '  +0.05%  '
$ws.Range("E37").Value2 = '  +3.64%  '
$ws.Range("D38").Value2 = '16.82'
$ws.Range("E38").Value2 = '  +8.12%  '
$ws.Range("B39").Value2 = 'Kaspa'
$ws.Range("C39").Value2 = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D39").Value2 = '0.105'
$ws.Range("E39").Value2 = '  +5.17%  '
$ws.Range("B40").Value2 = 'Stellar'
$ws.Range("C40").Value2 = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D40").Value2 = '0.116'
$ws.Range("E40").Value2 = '  +3.11%  '
$ws.Range("D41").Value2 = '1.83'
$ws.Range("E41").Value2 = '  +3.83%  '
$ws.Range("D42").Value2 = '3.94'
$ws.Range("E42").Value2 = '  +5.97%  '
$ws.Range("D43").Value2 = '2.078.99'
$ws.Range("E43").Value2 = '  -0.89%  '
$ws.Range("D44").Value2 = '19.66'
$ws.Range("E44").Value2 = '  +4.90%  '
$ws.Range("B45").Value2 = 'FraxShare'
$ws.Range("C45").Value2 = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D45").Value2 = '10.44'
$ws.Range("E45").Value2 = '  +2.89%  '
$ws.Range("B46").Value2 = 'VeChain'
$ws.Range("C46").Value2 = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D46").Value2 = '0.0280'
$ws.Range("E46").Value2 = '  +3.33%  '
$ws.Range("E47").Value2 = '  +7.16%  '
$ws.Range("E48").Value2 = '  +4.43%  '
$ws.Range("D49").Value2 = '1.54'
$ws.Range("E49").Value2 = '  +3.75%  '
$ws.Range("D50").Value2 = '73.19'
$ws.Range("E50").Value2 = '  +7.76%  '
$ws.Range("E51").Value2 = '  +3.52%  '

Write-Host "Applied" 97 "cell updates."
